# WAZIR_JULY_ATTENDANCE_POWER SYSTEM PROTECTION_EED_2024.xlsx
# The attendance workbook was rolled over from MAY to JULY: the sheet tab is
# renamed, the sheet-scoped Print_Area defined name is repointed at the new
# tab name, the July-21st attendance cell is filled in (present for the day),
# and the on-screen scroll/selection state is moved down to where the user
# left off editing.

$wb = $excel.ActiveWorkbook

# Rename the attendance sheet MAY_2024 -> JULY_2024
$ws = $wb.Worksheets.Item("MAY_2024")
$ws.Name = "JULY_2024"
$ws.Activate()

# The sheet-local "Print_Area" defined name still points at the old sheet
# name after the rename (the rename only updates the Name's own prefix) -
# repoint its formula at the renamed sheet.
$printArea = $wb.Names.Item("JULY_2024!Print_Area")
$printArea.RefersTo = "=JULY_2024!`$A`$1:`$H`$26"

# Mark attendance for row 21 (student 23ME09) on the "D" class day - the
# dependent SUM (column G) and percentage (column H) formulas recalc
# automatically.
$ws.Range("D21").Value = 3

# Move the viewport/selection to where editing left off.
$ws.Range("E24").Select()

$excel.Calculate()
